# "added the lift slide"
#
# Adds a new text box (bullet-style notes about the Lift class variables
# and methods) to slide 3 ("The Lift Class"), to the left of the existing
# elevator picture.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# The target shape is id=4 / name="TextBox 3" even though this slide only
# has 2 shapes (Title 1, Picture 2) right now. PowerPoint never re-uses
# shape ids, so an earlier shape must have been created (and removed)
# before this one during the original editing session. Reproduce that so
# the id/name numbering of the final textbox lines up.
$placeholderShape = $s.Shapes.AddTextbox(1, 0, 0, 10, 10)
$placeholderShape.Delete()

# Exact EMU target geometry (off x/y, ext cx/cy), expressed in points.
$left   = 26.181811023622046
$top    = 174.7636220472441
$width  = 348.218188976378
$height = 378.05622047244094

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "TextBox 3"
$tb.Fill.Visible = 0

$tf = $tb.TextFrame
$tf.WordWrap = -1

$tr = $tf.TextRange
$tr.Text = "Variables:"
$tr.InsertAfter("`r")
$tb.TextFrame.TextRange.Paragraphs(2).Text = "-A move variable that controls the movement of the elevator."
$tb.TextFrame.TextRange.InsertAfter("`r")
$tb.TextFrame.TextRange.Paragraphs(3).Text = "-A door variable that is either true for open or false for closed."
$tb.TextFrame.TextRange.InsertAfter("`r")
$tb.TextFrame.TextRange.Paragraphs(4).Text = "-A floor variable that tells the user the floor there going too."
$tb.TextFrame.TextRange.InsertAfter("`r")
$tb.TextFrame.TextRange.Paragraphs(5).Text = "-A position variable that shows the position of the elevator/where they currently are."
$tb.TextFrame.TextRange.InsertAfter("`r")
$tb.TextFrame.TextRange.Paragraphs(6).Text = "-Lift is the constructor."
$tb.TextFrame.TextRange.InsertAfter("`r")
$tb.TextFrame.TextRange.Paragraphs(7).Text = "-The "
$tb.TextFrame.TextRange.Paragraphs(7).InsertAfter("Move() method calculates the movement and the Door() method controls the doors.")
$tb.TextFrame.TextRange.InsertAfter("`r")
$tb.TextFrame.TextRange.Paragraphs(8).Text = ""
$tb.TextFrame.TextRange.InsertAfter("`r")
$tb.TextFrame.TextRange.Paragraphs(9).Text = ""
$tb.TextFrame.TextRange.InsertAfter("`r")
$tb.TextFrame.TextRange.Paragraphs(10).Text = ""

# Turning autosize on recalculates Height to fit the current text, so pin
# the geometry back to the authored size (the tiny +0.00001pt nudge on
# Height compensates for float rounding on the way to/from EMU).
$tf.AutoSize = 1
$tb.Height = 378.05623
